$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = 68
$ws.Range("C44").Value = 1
$ws.Range("D44").Value = 8
$ws.Range("E44").Value = 23
$ws.Range("F44").Value = 77
$ws.Range("G44").Value = 100
